$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

# A new "EditDescription" column is being introduced to the left of the
# existing "EditCategory" column (old column E), pushing it and everything
# to its right one column over (E:M -> F:N).
$ws.Range("E1").EntireColumn.Insert()

# Match the new column's width to its neighbour (old column E / new column F)
# so the pair renders the same as before the insert.
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(6).ColumnWidth()

# Fill in the new "Edit Description" column (header + sample data first,
# matching the order the shared strings show up in the saved file).
$ws.Range("E1").Value = "EditDescription"
$ws.Range("E2").Value = "Edit_Description"

# Refresh the sample category/name values to the new GNB705 test data.
$ws.Range("B2").Value = "AutomationTestGNB705"
$ws.Range("C2").Value = "AutomationTestGNB705"

# New value for the (now shifted) EditCategory sample column.
$ws.Range("F2").Value = "Edit_Automation_GNB705"

# Move the active selection to C10, matching the author's last selection.
[void]$ws.Range("C10").Select()
